# Issue #38 plus on add image is cross
# Add three new rows (36, 37, 38) to the "Issues" log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# Row 36: REST server check
$ws.Range("A36").Value = 36
$ws.Range("B36").Value = 1
$ws.Range("E36").Value = "REST server check"
$ws.Range("H36").Value = "Check connectivity to Server on settings page"

# Row 37: REST ping service
$ws.Range("A37").Value = 37
$ws.Range("B37").Value = 1
$ws.Range("E37").Value = "REST ping service"
$ws.Range("F37").Value = 37
$ws.Range("H37").Value = "Add a rest ping service"

# Row 38: plus on add image is cross
$ws.Range("A38").Value = 38
$ws.Range("B38").Value = 1
$ws.Range("D38").Value = "Bug"
$ws.Range("E38").Value = "plus on add image is cross"

# Leave selection on the last entered cell, matching the authored edit
$ws.Range("E38").Select()
